$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Events")

# New text for B65 (e052 "Pivot Tank") - shortened, no image, new wording
$e052Text = @'
<Bold>e052 Pivot Tank</Bold> 
<InlineUIContainer><Button Content='r4.74.1' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r8.46' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
If you want the Sherman to face a different sector, select buttons here:
<LineBreak/><LineBreak/>
                                           <InlineUIContainer><Button Content='   -   ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>      
<InlineUIContainer><Button Content='   +   ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$e052Text = $e052Text.TrimEnd("`r","`n")

# New text for B66 (e052a "Pivot Turret") - shortened, no image, new wording
$e052aText = @'
<Bold>e052a Pivot Turret</Bold> 
<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r8.24' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
If you want the turret to face a different sector, select buttons here: 
<LineBreak/><LineBreak/>
                                             <InlineUIContainer><Button Content='  -  ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>           
<InlineUIContainer><Button Content='  +  ' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>
<LineBreak/><LineBreak/>
'@
$e052aText = $e052aText.TrimEnd("`r","`n")

# New (truncated) text for B69 (e053b "Main Gun Firing - Target Selected")
# The paragraph about AAR / Direct / Area fire was removed.
$e053bText = @'
<Bold>e053b Main Gun Firing - Target Selected</Bold> 
<InlineUIContainer><Button Content='r4.74.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>   
<InlineUIContainer><Button Content='r9.0' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> 
<LineBreak/><LineBreak/>
'@
$e053bText = $e053bText.TrimEnd("`r","`n")

$ws.Range("B65").Value = $e052Text
$ws.Range("B66").Value = $e052aText
$ws.Range("B69").Value = $e053bText

# Row heights shrank now that the text is shorter
$ws.Rows.Item(65).RowHeight = 135
$ws.Rows.Item(66).RowHeight = 135
$ws.Rows.Item(69).RowHeight = 60

# Update the saved selection to match (scrolled down a bit further, selection on B65)
$ws.Activate()
$ws.Range("B65").Select()
